$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I13").Value = '%'
$ws.Range("J13").Value = 'Uninterpretable'
$ws.Range("I17").Value = 'sv'
$ws.Range("J17").Value = 'Statement-opinion'
$ws.Range("I33").Value = 'sd'
$ws.Range("J33").Value = 'Statement-non-opinion'
$ws.Range("I37").Value = 'b'
$ws.Range("J37").Value = 'Acknowledge (Backchannel)'
$ws.Range("I52").Value = 'sd'
$ws.Range("J52").Value = 'Statement-non-opinion'
$ws.Range("I65").Value = 'sv'
$ws.Range("J65").Value = 'Statement-opinion'
$ws.Range("I69").Value = 'sd'
$ws.Range("J69").Value = 'Statement-non-opinion'
$ws.Range("I83").Value = 'sv'
$ws.Range("J83").Value = 'Statement-opinion'
$ws.Range("I95").Value = 'sd'
$ws.Range("J95").Value = 'Statement-non-opinion'
$ws.Range("I100").Value = 'sd'
$ws.Range("J100").Value = 'Statement-non-opinion'
$ws.Range("I112").Value = 'sv'
$ws.Range("J112").Value = 'Statement-opinion'
$ws.Range("I114").Value = 'sv'
$ws.Range("J114").Value = 'Statement-opinion'
$ws.Range("I120").Value = 'ba'
$ws.Range("J120").Value = 'Appreciation'
$ws.Range("I126").Value = 'sd'
$ws.Range("J126").Value = 'Statement-non-opinion'
$ws.Range("I130").Value = 'aa'
$ws.Range("J130").Value = 'Agree/Accept'
$ws.Range("I137").Value = 'ba'
$ws.Range("J137").Value = 'Appreciation'
$ws.Range("I168").Value = 'qy'
$ws.Range("J168").Value = 'Yes-No-Question'
$ws.Range("I180").Value = 'sv'
$ws.Range("J180").Value = 'Statement-opinion'
$ws.Range("I191").Value = '%'
$ws.Range("J191").Value = 'Uninterpretable'
$ws.Range("I195").Value = 'aa'
$ws.Range("J195").Value = 'Agree/Accept'
$ws.Range("I216").Value = 'b'
$ws.Range("J216").Value = 'Acknowledge (Backchannel)'
$ws.Range("I221").Value = 'aa'
$ws.Range("J221").Value = 'Agree/Accept'
$ws.Range("I222").Value = 'b'
$ws.Range("J222").Value = 'Acknowledge (Backchannel)'
$ws.Range("I230").Value = 'sd'
$ws.Range("J230").Value = 'Statement-non-opinion'
$ws.Range("I233").Value = 'aa'
$ws.Range("J233").Value = 'Agree/Accept'
$ws.Range("I237").Value = 'aa'
$ws.Range("J237").Value = 'Agree/Accept'
$ws.Range("I248").Value = 'aa'
$ws.Range("J248").Value = 'Agree/Accept'
$ws.Range("I256").Value = 'b'
$ws.Range("J256").Value = 'Acknowledge (Backchannel)'
$ws.Range("I266").Value = 'b'
$ws.Range("J266").Value = 'Acknowledge (Backchannel)'
$ws.Range("I270").Value = 'sd'
$ws.Range("J270").Value = 'Statement-non-opinion'
$ws.Range("I273").Value = 'b'
$ws.Range("J273").Value = 'Acknowledge (Backchannel)'
$ws.Range("I289").Value = 'b'
$ws.Range("J289").Value = 'Acknowledge (Backchannel)'
$ws.Range("I290").Value = 'sd'
$ws.Range("J290").Value = 'Statement-non-opinion'
$ws.Range("I293").Value = 'sd'
$ws.Range("J293").Value = 'Statement-non-opinion'
$ws.Range("I297").Value = 'sd'
$ws.Range("J297").Value = 'Statement-non-opinion'
$ws.Range("I307").Value = 'sv'
$ws.Range("J307").Value = 'Statement-opinion'
$ws.Range("I311").Value = '%'
$ws.Range("J311").Value = 'Uninterpretable'
$ws.Range("I343").Value = 'ba'
$ws.Range("J343").Value = 'Appreciation'
$ws.Range("I347").Value = 'sd'
$ws.Range("J347").Value = 'Statement-non-opinion'
